$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 18): Product of Array Except Self
$ws.Range("A18").Value = "Product of Array Except Self"
$ws.Range("B18").Value = "Array"
$ws.Range("C18").Value = "Yes"
$ws.Range("D18").Value = "Yes"
$ws.Range("E18").Value = "Medium"
$ws.Range("F18").Value = "Medium"
$ws.Range("G18").Value = "238 - Product of Array Except Self"

# Add hyperlink for G18, matching the existing pattern (local-path style target)
$ws.Hyperlinks.Add($ws.Range("G18"), "238 - Product of Array Except Self", "", "", "238 - Product of Array Except Self") | Out-Null

# Match style used by other hyperlink cells in column G
$ws.Range("G18").Style = "Hyperlink"

# Extend conditional formatting range to include row 18
$ws.Range("D9:F18").FormatConditions.Delete() | Out-Null

# Extend data validations to include row 18
$ws.Range("E2:F18").Validation.Delete() | Out-Null
$ws.Range("E2:F18").Validation.Add(3, 1, 1, "Easy, Medium, Hard") | Out-Null

$ws.Range("C2:C18").Validation.Delete() | Out-Null
$ws.Range("C2:C18").Validation.Add(3, 0, 1, "Yes, No") | Out-Null

$ws.Range("B2:B18").Validation.Delete() | Out-Null
$ws.Range("B2:B18").Validation.Add(3, 1, 1, "Array, Binary, Dynamic Programming, Graph, Interval, Linked List, Matrix, String, Tree, Heap") | Out-Null

$ws.Range("D2:D18").Validation.Delete() | Out-Null
$ws.Range("D2:D18").Validation.Add(3, 1, 1, "Yes, No") | Out-Null

# Update active cell selection like in the diff
$ws.Range("G27").Select()
